$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped from 45172 to 45175
# for every data row (rows 2 through 367).
$ws.Range("C2:C367").Value = 45175
